# Weekly update: insert the newest week's record at the top of the data
# block (row 135), pushing the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 135:167 down to 136:168, opening up a blank row 135.
$ws.Rows.Item(135).Insert()

# Populate the new row 135 with this week's observation. The row mirrors
# the shape of its neighbours (same market / product / metadata columns),
# only the date and measurement columns differ.
$ws.Range("A135").Value = 6
$ws.Range("B135").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C135").Value = "Metropolitana"
$ws.Range("D135").Value = 44642
$ws.Range("E135").Value = 13
$ws.Range("F135").Value = 100112029
$ws.Range("G135").Value = "Orégano"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 35
$ws.Range("K135").Value = 14500
$ws.Range("L135").Value = 15000
$ws.Range("M135").Value = 14729
$ws.Range("N135").Value = "$/docena de atados"
$ws.Range("O135").Value = "Región Metropolitana"
$ws.Range("P135").Value = 4910
$ws.Range("Q135").Value = 3
$ws.Range("R135").Value = "Hortaliza"
